$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the description text for the "Most popular properties" checklist item
$ws.Range("B22").Value = "Most popular properties on gallery page. Similar Houses based on location on single page."

# Mark checklist item A5 as done (TRUE) - Google Maps plotting feature
$ws.Range("A5").Value = $true

# Update selection / view state to match the diff
$ws.Range("B28").Select()
$excel.ActiveWindow.ScrollRow = 16
